# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit refresh to Sheets/Excalibur_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1447.0869
$ws.Range("J17").Value = 1518.9
$ws.Range("L17").Value = 4556.700000000001
$ws.Range("N17").Value = -4892.700000000001
$ws.Range("H112").Value = 958.4138
$ws.Range("J112").Value = 933.1111
$ws.Range("L112").Value = 2799.3333
$ws.Range("N112").Value = -5015.3333
$ws.Range("H137").Value = 26318072
$ws.Range("I137").Value = 55557930
$ws.Range("J137").Value = 2201.3
$ws.Range("K137").Value = 166673790
$ws.Range("L137").Value = 6603.900000000001
$ws.Range("M137").Value = -166671240
$ws.Range("N137").Value = -11703.9
$ws.Range("H138").Value = 749.4286
$ws.Range("I138").Value = 749.4286
$ws.Range("K138").Value = 2248.2858
$ws.Range("M138").Value = 2891.7142
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2448.353
$ws.Range("I45").Value = 2187.5715
$ws.Range("K45").Value = 2187.5715
$ws.Range("M45").Value = -1810.5715
$ws.Range("H61").Value = 1960.081
$ws.Range("I61").Value = 1763.5161
$ws.Range("K61").Value = 1763.5161
$ws.Range("M61").Value = -1551.5161
$ws.Range("H74").Value = 1448.2954
$ws.Range("I74").Value = 744.75757
$ws.Range("K74").Value = 744.75757
$ws.Range("M74").Value = 129.24243
$ws.Range("H77").Value = 1448.2954
$ws.Range("I77").Value = 744.75757
$ws.Range("K77").Value = 3723.78785
$ws.Range("M77").Value = 644.2121500000003
$ws.Range("H136").Value = 1960.081
$ws.Range("I136").Value = 1763.5161
$ws.Range("K136").Value = 5290.5483
$ws.Range("M136").Value = -2740.5483
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2318.6667
$ws.Range("I20").Value = 2246.2222
$ws.Range("J20").Value = 2536
$ws.Range("K20").Value = 2246.2222
$ws.Range("L20").Value = 2536
$ws.Range("M20").Value = -1999.2222
$ws.Range("N20").Value = -3030
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5527.1704
$ws.Range("I31").Value = 5622.273
$ws.Range("J31").Value = 5443.48
$ws.Range("K31").Value = 5622.273
$ws.Range("L31").Value = 5443.48
$ws.Range("M31").Value = -5327.273
$ws.Range("N31").Value = -6033.48
$ws.Range("H34").Value = 5527.1704
$ws.Range("I34").Value = 5622.273
$ws.Range("J34").Value = 5443.48
$ws.Range("K34").Value = 5622.273
$ws.Range("L34").Value = 5443.48
$ws.Range("M34").Value = -5420.273
$ws.Range("N34").Value = -5847.48
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("K39").Value = 1000
$ws.Range("M39").Value = -609
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("K49").Value = 1000
$ws.Range("M49").Value = -818
$ws.Range("H58").Value = 2013.5952
$ws.Range("I58").Value = 1256.6774
$ws.Range("J58").Value = 4146.727
$ws.Range("K58").Value = 1256.6774
$ws.Range("L58").Value = 4146.727
$ws.Range("M58").Value = -1053.6774
$ws.Range("N58").Value = -4552.727
$ws.Range("H136").Value = 2013.5952
$ws.Range("I136").Value = 1256.6774
$ws.Range("J136").Value = 4146.727
$ws.Range("K136").Value = 3770.0322
$ws.Range("L136").Value = 12440.181
$ws.Range("M136").Value = -1220.0322
$ws.Range("N136").Value = -17540.181
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 952.4286
$ws.Range("I5").Value = 744.5
$ws.Range("J5").Value = 1108.375
$ws.Range("K5").Value = 2233.5
$ws.Range("L5").Value = 3325.125
$ws.Range("M5").Value = -2121.5
$ws.Range("N5").Value = -3549.125
$ws.Range("H46").Value = 200
$ws.Range("I46").Value = 200
$ws.Range("K46").Value = 600
$ws.Range("M46").Value = -509
$ws.Range("H50").Value = 680.125
$ws.Range("I50").Value = 525.7273
$ws.Range("K50").Value = 1577.1819
$ws.Range("M50").Value = -1096.1819
$ws.Range("H53").Value = 680.125
$ws.Range("I53").Value = 525.7273
$ws.Range("K53").Value = 1577.1819
$ws.Range("M53").Value = -1096.1819
$ws.Range("H119").Value = 2714.2856
$ws.Range("I119").Value = 2166.6667
$ws.Range("J119").Value = 3125
$ws.Range("K119").Value = 6500.000100000001
$ws.Range("L119").Value = 9375
$ws.Range("M119").Value = -1662.000100000001
$ws.Range("N119").Value = -19051
$ws.Range("H120").Value = 19666.545
$ws.Range("I120").Value = 8028.5713
$ws.Range("J120").Value = 40033
$ws.Range("K120").Value = 24085.7139
$ws.Range("L120").Value = 120099
$ws.Range("M120").Value = -19247.7139
$ws.Range("N120").Value = -129775
$ws.Range("H130").Value = 1823
$ws.Range("I130").Value = 1906.25
$ws.Range("J130").Value = 1490
$ws.Range("K130").Value = 5718.75
$ws.Range("L130").Value = 4470
$ws.Range("M130").Value = -698.75
$ws.Range("N130").Value = -14510
$ws.Range("H135").Value = 952.4286
$ws.Range("I135").Value = 744.5
$ws.Range("J135").Value = 1108.375
$ws.Range("K135").Value = 6700.5
$ws.Range("L135").Value = 9975.375
$ws.Range("M135").Value = -4165.5
$ws.Range("N135").Value = -15045.375
$ws.Range("H139").Value = 61471.41
$ws.Range("I139").Value = 69108
$ws.Range("K139").Value = 207324
$ws.Range("M139").Value = -202184
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1181.2222
$ws.Range("I82").Value = 870
$ws.Range("J82").Value = 1803.6666
$ws.Range("K82").Value = 870
$ws.Range("L82").Value = 1803.6666
$ws.Range("M82").Value = -509
$ws.Range("N82").Value = -2525.6666
$ws.Range("H85").Value = 1181.2222
$ws.Range("I85").Value = 870
$ws.Range("J85").Value = 1803.6666
$ws.Range("K85").Value = 870
$ws.Range("L85").Value = 1803.6666
$ws.Range("M85").Value = 378
$ws.Range("N85").Value = -4299.6666
$ws.Range("H93").Value = 1941.421
$ws.Range("I93").Value = 1911.6875
$ws.Range("J93").Value = 2100
$ws.Range("K93").Value = 1911.6875
$ws.Range("L93").Value = 2100
$ws.Range("M93").Value = -663.6875
$ws.Range("N93").Value = -4596
$ws.Range("H100").Value = 10646.308
$ws.Range("I100").Value = 3199.8333
$ws.Range("K100").Value = 3199.8333
$ws.Range("M100").Value = -2658.8333
$ws.Range("H131").Value = 18642.8
$ws.Range("I131").Value = 10296
$ws.Range("J131").Value = 31163
$ws.Range("K131").Value = 10296
$ws.Range("L131").Value = 31163
$ws.Range("M131").Value = -5256
$ws.Range("N131").Value = -41243
$ws.Range("H136").Value = 4109.2
$ws.Range("I136").Value = 2636.5
$ws.Range("K136").Value = 7909.5
$ws.Range("M136").Value = -5359.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 44165.832
$ws.Range("J63").Value = 50199
$ws.Range("L63").Value = 50199
$ws.Range("N63").Value = -51447
$ws.Range("H66").Value = 44165.832
$ws.Range("J66").Value = 50199
$ws.Range("L66").Value = 150597
$ws.Range("N66").Value = -156837
$ws.Range("H68").Value = 66168.625
$ws.Range("I68").Value = 43089.668
$ws.Range("J68").Value = 80016
$ws.Range("K68").Value = 43089.668
$ws.Range("L68").Value = 80016
$ws.Range("M68").Value = -42278.668
$ws.Range("N68").Value = -81638
$ws.Range("H71").Value = 66168.625
$ws.Range("I71").Value = 43089.668
$ws.Range("J71").Value = 80016
$ws.Range("K71").Value = 129269.004
$ws.Range("L71").Value = 240048
$ws.Range("M71").Value = -125213.004
$ws.Range("N71").Value = -248160
$ws.Range("H116").Value = 200644
$ws.Range("J116").Value = 200644
$ws.Range("L116").Value = 200644
$ws.Range("N116").Value = -209822
$ws.Range("H136").Value = 5954908
$ws.Range("I136").Value = 6175089.5
$ws.Range("K136").Value = 18525268.5
$ws.Range("M136").Value = -18522718.5
